$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the values of a handful of columns between row 2 and row 3
# (A, I, Q, R, AC) - effectively exchanging the two records' Id/Antal/Ost/Nord/
# Publik kommentar values while leaving the rest of each row untouched.

$a2  = $ws.Range("A2").Value2
$i2  = $ws.Range("I2").Value2
$q2  = $ws.Range("Q2").Value2
$r2  = $ws.Range("R2").Value2
$ac2 = $ws.Range("AC2").Value2

$a3  = $ws.Range("A3").Value2
$i3  = $ws.Range("I3").Value2
$q3  = $ws.Range("Q3").Value2
$r3  = $ws.Range("R3").Value2
$ac3 = $ws.Range("AC3").Value2

$ws.Range("A2").Value = $a3
$ws.Range("I2").Value = $i3
$ws.Range("Q2").Value = $q3
$ws.Range("R2").Value = $r3
$ws.Range("AC2").Value = $ac3

$ws.Range("A3").Value = $a2
$ws.Range("I3").Value = $i2
$ws.Range("Q3").Value = $q2
$ws.Range("R3").Value = $r2
$ws.Range("AC3").Value = $ac2
